$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the image path value used in rows 5 and 6 (column E)
$ws.Range("E5").Value = "C:\Users\username\Desktop\Untitled.png"
$ws.Range("E6").Value = "C:\Users\username\Desktop\Untitled.png"

# Update the active selection to match the new state
$ws.Range("E7").Select()
